$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.300.32'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '3.529.09'
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''597.39'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").Value = '''171.94'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +8.01%  '
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").Value = '''0.439'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '4.131.88'
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = '''28.67'
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("D16").Value = '67.257.72'
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").Value = '3.513.82'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '''6.38'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").Value = '''14.21'
$ws.Range("E19").Value = '  +1.28%  '
$ws.Range("D20").Value = '''397.71'
$ws.Range("E20").Value = '  +2.38%  '
$ws.Range("D21").Value = '''8.04'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").Value = '''0.541'
$ws.Range("E23").Value = '  +2.50%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '''0.0000124'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '''10.27'
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '''6.32'
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").Value = '''1.48'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("D32").Value = '''24.19'
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("D33").Value = '''7.43'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("E34").Value = '  +4.81%  '
$ws.Range("D35").Value = '''163.88'
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("D36").Value = '''0.899'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("E38").Value = '  +3.62%  '
$ws.Range("D39").Value = '''6.87'
$ws.Range("E39").Value = '  +2.25%  '
$ws.Range("D40").Value = '''0.0751'
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").Value = '''26.73'
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("D42").Value = '''27.06'
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''2.62'
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.811.70'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '''42.94'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '''343.29'
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").Value = '''33.60'
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").Value = '''6.57'
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").Value = '''0.859'
$ws.Range("E51").Value = '  +0.91%  '
